$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "317.61" }
    @{ Cell = "E2"; Value = "3.78%" }
    @{ Cell = "D3"; Value = "39.80" }
    @{ Cell = "E3"; Value = "1.86%" }
    @{ Cell = "D4"; Value = "5.156" }
    @{ Cell = "E4"; Value = "0.87%" }
    @{ Cell = "D5"; Value = "0.08236" }
    @{ Cell = "D6"; Value = "2.053" }
    @{ Cell = "E6"; Value = "6.22%" }
    @{ Cell = "D7"; Value = "8.387" }
    @{ Cell = "E7"; Value = "4.09%" }
    @{ Cell = "D8"; Value = "0.9405" }
    @{ Cell = "E8"; Value = "1.57%" }
    @{ Cell = "D9"; Value = "0.1353" }
    @{ Cell = "E9"; Value = "-2.87%" }
    @{ Cell = "D10"; Value = "0.1995" }
    @{ Cell = "E10"; Value = "3.97%" }
    @{ Cell = "D11"; Value = "0.09086" }
    @{ Cell = "E11"; Value = "0.66%" }
    @{ Cell = "D12"; Value = "0.03514" }
    @{ Cell = "E12"; Value = "-0.07%" }
    @{ Cell = "D13"; Value = "0.09805" }
    @{ Cell = "E13"; Value = "0.12%" }
    @{ Cell = "D14"; Value = "0.001411" }
    @{ Cell = "E14"; Value = "1.33%" }
    @{ Cell = "D15"; Value = "0.006115" }
    @{ Cell = "E15"; Value = "4.99%" }
    @{ Cell = "D16"; Value = "3.687" }
    @{ Cell = "E16"; Value = "-2.06%" }
    @{ Cell = "D17"; Value = "4.324" }
    @{ Cell = "E17"; Value = "2.98%" }
    @{ Cell = "D18"; Value = "3.349" }
    @{ Cell = "E18"; Value = "-0.92%" }
    @{ Cell = "D19"; Value = "0.3495" }
    @{ Cell = "E19"; Value = "0.96%" }
    @{ Cell = "E20"; Value = "-0.15%" }
    @{ Cell = "D21"; Value = "4.960" }
    @{ Cell = "E21"; Value = "6.04%" }
    @{ Cell = "D22"; Value = "0.2454" }
    @{ Cell = "E22"; Value = "1.59%" }
    @{ Cell = "D23"; Value = "0.04345" }
    @{ Cell = "E23"; Value = "-0.70%" }
    @{ Cell = "D24"; Value = "0.001240" }
    @{ Cell = "E24"; Value = "2.81%" }
    @{ Cell = "D25"; Value = "0.004798" }
    @{ Cell = "E25"; Value = "12.16%" }
    @{ Cell = "E26"; Value = "-0.01%" }
    @{ Cell = "E27"; Value = "-10.03%" }
    @{ Cell = "E39"; Value = "12.71%" }
    @{ Cell = "D40"; Value = "0.05181" }
    @{ Cell = "E40"; Value = "3.02%" }
    @{ Cell = "D41"; Value = "0.007768" }
    @{ Cell = "E41"; Value = "3.19%" }
    @{ Cell = "D42"; Value = "0.01047" }
    @{ Cell = "E42"; Value = "7.80%" }
    @{ Cell = "D43"; Value = "0.1408" }
    @{ Cell = "E43"; Value = "4.69%" }
    @{ Cell = "D44"; Value = "0.002045" }
    @{ Cell = "E44"; Value = "-2.27%" }
    @{ Cell = "D45"; Value = "0.009305" }
    @{ Cell = "E45"; Value = "-5.04%" }
    @{ Cell = "D46"; Value = "0.00006610" }
    @{ Cell = "E46"; Value = "6.24%" }
    @{ Cell = "E47"; Value = "-0.11%" }
    @{ Cell = "D48"; Value = "0.002896" }
    @{ Cell = "E48"; Value = "4.06%" }
    @{ Cell = "E49"; Value = "-6.25%" }
    @{ Cell = "E50"; Value = "-0.11%" }
    @{ Cell = "E51"; Value = "-0.11%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
